# Update Belias_Profits market-price / leve-profit figures per scheduled Universalis data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 133.90909
$ws.Range("I9").Value = 190.6
$ws.Range("K9").Value = 190.6
$ws.Range("M9").Value = -21.59999999999999

$ws.Range("H53").Value = 907.6957
$ws.Range("I53").Value = 1233.1875
$ws.Range("K53").Value = 1233.1875
$ws.Range("M53").Value = -596.1875

$ws.Range("H64").Value = 3242
$ws.Range("I64").Value = 2873.75
$ws.Range("J64").Value = 3340.2
$ws.Range("K64").Value = 2873.75
$ws.Range("L64").Value = 3340.2
$ws.Range("M64").Value = -2625.75
$ws.Range("N64").Value = -3836.2

$ws.Range("H67").Value = 3242
$ws.Range("I67").Value = 2873.75
$ws.Range("J67").Value = 3340.2
$ws.Range("K67").Value = 2873.75
$ws.Range("L67").Value = 3340.2
$ws.Range("M67").Value = -2015.75
$ws.Range("N67").Value = -5056.2

$ws.Range("H103").Value = 5904.1
$ws.Range("J103").Value = 879.625
$ws.Range("L103").Value = 2638.875
$ws.Range("N103").Value = -3810.875

$ws.Range("H106").Value = 5308.5
$ws.Range("I106").Value = 5300.5557
$ws.Range("J106").Value = 5311.905
$ws.Range("K106").Value = 5300.5557
$ws.Range("L106").Value = 5311.905
$ws.Range("M106").Value = -4669.5557
$ws.Range("N106").Value = -6573.905

$ws.Range("H116").Value = 2789.4443
$ws.Range("I116").Value = 2684.1667
$ws.Range("K116").Value = 2684.1667
$ws.Range("M116").Value = 757.8332999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2800
$ws.Range("I102").Value = 2650
$ws.Range("J102").Value = 2900
$ws.Range("K102").Value = 2650
$ws.Range("L102").Value = 2900
$ws.Range("M102").Value = -1028
$ws.Range("N102").Value = -6144

$ws.Range("H123").Value = 21111.111
$ws.Range("J123").Value = 21111.111
$ws.Range("L123").Value = 21111.111
$ws.Range("N123").Value = -30911.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29259.334
$ws.Range("J55").Value = 29259.334
$ws.Range("L55").Value = 29259.334
$ws.Range("N55").Value = -29805.334

$ws.Range("H86").Value = 2341.3635
$ws.Range("I86").Value = 1778.7333
$ws.Range("J86").Value = 3547
$ws.Range("K86").Value = 1778.7333
$ws.Range("L86").Value = 3547
$ws.Range("M86").Value = -655.7333000000001
$ws.Range("N86").Value = -5793

$ws.Range("H89").Value = 2341.3635
$ws.Range("I89").Value = 1778.7333
$ws.Range("J89").Value = 3547
$ws.Range("K89").Value = 8893.666500000001
$ws.Range("L89").Value = 17735
$ws.Range("M89").Value = -3277.666500000001
$ws.Range("N89").Value = -28967

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 950
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 950
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 950
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1524

$ws.Range("H62").Value = 2338.0605
$ws.Range("I62").Value = 2166.6667
$ws.Range("J62").Value = 2638
$ws.Range("K62").Value = 2166.6667
$ws.Range("L62").Value = 2638
$ws.Range("M62").Value = -1542.6667
$ws.Range("N62").Value = -3886

$ws.Range("H65").Value = 2338.0605
$ws.Range("I65").Value = 2166.6667
$ws.Range("J65").Value = 2638
$ws.Range("K65").Value = 10833.3335
$ws.Range("L65").Value = 13190
$ws.Range("M65").Value = -7713.333500000001
$ws.Range("N65").Value = -19430

$ws.Range("H113").Value = 950
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 950
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5290

$ws.Range("H141").Value = 36358.273
$ws.Range("J141").Value = 36358.273
$ws.Range("L141").Value = 36358.273
$ws.Range("N141").Value = -46718.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 42127.418
$ws.Range("I2").Value = 25.1
$ws.Range("J2").Value = 72200.5
$ws.Range("K2").Value = 150.6
$ws.Range("L2").Value = 433203
$ws.Range("M2").Value = -37.60000000000002
$ws.Range("N2").Value = -433429

$ws.Range("H4").Value = 955.1724
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 1126.0869
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 3378.2607
$ws.Range("M4").Value = -788
$ws.Range("N4").Value = -3602.2607

$ws.Range("H7").Value = 16666699
$ws.Range("J7").Value = 66666700
$ws.Range("L7").Value = 200000100
$ws.Range("N7").Value = -200000324

$ws.Range("H10").Value = 134.5
$ws.Range("I10").Value = 65.666664
$ws.Range("J10").Value = 258.4
$ws.Range("K10").Value = 196.999992
$ws.Range("L10").Value = 775.1999999999999
$ws.Range("M10").Value = -57.99999199999999
$ws.Range("N10").Value = -1053.2

$ws.Range("H11").Value = 448.44446
$ws.Range("I11").Value = 300
$ws.Range("J11").Value = 467
$ws.Range("K11").Value = 900
$ws.Range("L11").Value = 1401
$ws.Range("M11").Value = -760
$ws.Range("N11").Value = -1681

$ws.Range("H12").Value = 315.35294
$ws.Range("I12").Value = 340.14285
$ws.Range("J12").Value = 298
$ws.Range("K12").Value = 1020.42855
$ws.Range("L12").Value = 894
$ws.Range("M12").Value = -847.4285500000001
$ws.Range("N12").Value = -1240

$ws.Range("H131").Value = 849.45
$ws.Range("I131").Value = 405.58334
$ws.Range("J131").Value = 909.9773
$ws.Range("K131").Value = 1216.75002
$ws.Range("L131").Value = 2729.9319
$ws.Range("M131").Value = 3823.24998
$ws.Range("N131").Value = -12809.9319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 18836.182
$ws.Range("J57").Value = 18836.182
$ws.Range("L57").Value = 18836.182
$ws.Range("N57").Value = -20476.182

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 50290.285
$ws.Range("J25").Value = 57338.668
$ws.Range("L25").Value = 57338.668
$ws.Range("N25").Value = -57798.668

$ws.Range("H55").Value = 362.03705
$ws.Range("I55").Value = 129.22223
$ws.Range("J55").Value = 478.44446
$ws.Range("K55").Value = 129.22223
$ws.Range("L55").Value = 478.44446
$ws.Range("M55").Value = 43.77777
$ws.Range("N55").Value = -824.4444599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 24282.5
$ws.Range("I75").Value = 8000
$ws.Range("J75").Value = 29710
$ws.Range("K75").Value = 8000
$ws.Range("L75").Value = 29710
$ws.Range("M75").Value = -7064
$ws.Range("N75").Value = -31582

$ws.Range("H78").Value = 24282.5
$ws.Range("I78").Value = 8000
$ws.Range("J78").Value = 29710
$ws.Range("K78").Value = 24000
$ws.Range("L78").Value = 89130
$ws.Range("M78").Value = -19320
$ws.Range("N78").Value = -98490
